# Apply cryptos list update (prices, volumes, and a row swap for Toncoin/InjectiveProtocol)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.167.97"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.949.57"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "377.04"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.51"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.41"
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0853"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "3.411.21"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.11"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.61"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "2.965.22"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.00"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.05"
$ws.Range("E18").Value = "  +48.45%  "
$ws.Range("D19").Value = "51.137.32"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("E20").Value = "  -6.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.48"
$ws.Range("E21").Value = "  -3.54%  "
$ws.Range("D22").Value = "0.0₃0956"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "265.86"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.76"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("E25").Value = "  +7.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.22"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.59"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -3.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.66"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.04"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.86"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.05"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.49"
$ws.Range("E35").Value = "  -3.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0443"
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  +3.91%  "
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.40"
$ws.Range("E40").Value = "  -4.12%  "
$ws.Range("E41").Value = "  -2.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.49"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.64"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.38"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.38"
$ws.Range("E45").Value = "  +3.50%  "
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.273"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").Value = "1.991.68"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0326"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("E51").Value = "  +2.27%  "
